$d = $word.ActiveDocument

# Replace the placeholder ID text in the first paragraph (this also removes
# the trailing single-space run that followed it, since the search text
# includes that trailing space while the replacement does not).
$d.Content.Find.Execute(
    "**ID__AFFARS_5319_topic_12__ID** ", $true, $false, $false, $false, $false,
    $true, 1, $false, "**ID__AFFARS_5319_1305__ID**", 2)

# Update the first paragraph's formatting: increase the left indent from
# 120 twips (6pt) to 225 twips (11.25pt) and add a paragraph border whose
# lines all have a 5pt text-to-border spacing (matching the pattern already
# used elsewhere in the document).
$p = $d.Paragraphs.Item(1)
$p.Format.LeftIndent = 11.25
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
